# Lecture 11 on data privacy — trim slide deck and fix title text.
$p = $ppt.ActivePresentation

# 1) Merge the "Fall " / "2021" runs on slide 1 into a single run
#    "Fall 2021" while keeping the dirty="0" run's formatting.
$s1 = $p.Slides.Item(1)
$tr = $s1.Shapes.Item(2).TextFrame.TextRange
$para2 = $tr.Paragraphs(2)
$yearRun = $para2.Characters(6, 4)
$yearRun.Text = "Fall 2021"
$para2 = $tr.Paragraphs(2)
$prefixRun = $para2.Characters(1, 5)
$prefixRun.Text = ""

# 2) Drop the trailing 12 slides (old slides 27-38: browser/O.S./email
#    threat content) that aren't part of this lecture anymore.
for ($i = $p.Slides.Count; $i -ge 27; $i--) {
    $p.Slides.Item($i).Delete()
}
